$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61
$ws.Range("A61").Value = 60
$ws.Range("B61").Value = "29/12/2025 04:20"
$ws.Range("C61").Value = 400
$ws.Range("D61").Value = "Conhecimentos Específicos"
$ws.Range("E61").Value = "Gestão da Produção e Operações"
$ws.Range("F61").Value = 1

# Row 62
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = "29/12/2025 04:21"
$ws.Range("C62").Value = 491
$ws.Range("D62").Value = "Conhecimentos Específicos"
$ws.Range("E62").Value = "Processos Decisórios"
$ws.Range("F62").Value = 1

# Row 63
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = "29/12/2025 04:25"
$ws.Range("C63").Value = 536
$ws.Range("D63").Value = "Conhecimentos Específicos"
$ws.Range("E63").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("F63").Value = 1

# Row 64
$ws.Range("A64").Value = 63
$ws.Range("B64").Value = "29/12/2025 04:27"
$ws.Range("C64").Value = 545
$ws.Range("D64").Value = "Conhecimentos Específicos"
$ws.Range("E64").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("F64").Value = 0

# Row 65
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = "29/12/2025 04:28"
$ws.Range("C65").Value = 1538
$ws.Range("D65").Value = "Contabilidade Gerencial"
$ws.Range("E65").Value = "Alavancagem Financeira"
$ws.Range("F65").Value = 0

# Row 66
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = "29/12/2025 04:30"
$ws.Range("C66").Value = "'1568"
$ws.Range("D66").Value = "Contabilidade Gerencial"
$ws.Range("E66").Value = "EBITDA"
$ws.Range("F66").Value = 0

# Row 67
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = "29/12/2025 04:32"
$ws.Range("C67").Value = 1565
$ws.Range("D67").Value = "Contabilidade Gerencial"
$ws.Range("E67").Value = "EBITDA"
$ws.Range("F67").Value = 0

# Row 68
$ws.Range("A68").Value = 67
$ws.Range("B68").Value = "29/12/2025 04:34"
$ws.Range("C68").Value = 1523
$ws.Range("D68").Value = "Contabilidade Gerencial"
$ws.Range("E68").Value = "DRE"
$ws.Range("F68").Value = 0

# Row 69
$ws.Range("A69").Value = 68
$ws.Range("B69").Value = "29/12/2025 04:41"
$ws.Range("C69").Value = 1007
$ws.Range("D69").Value = "Estatística"
$ws.Range("E69").Value = "Medidas de Variabilidade"
$ws.Range("F69").Value = 0

# Row 70
$ws.Range("A70").Value = 69
$ws.Range("B70").Value = "29/12/2025 04:46"
$ws.Range("C70").Value = 1191
$ws.Range("D70").Value = "Estatística"
$ws.Range("E70").Value = "Testes de Hipóteses"
$ws.Range("F70").Value = 0

# Row 71
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = "29/12/2025 04:49"
$ws.Range("C71").Value = 1241
$ws.Range("D71").Value = "Estatística"
$ws.Range("E71").Value = "ANOVA"
$ws.Range("F71").Value = 0

# Row 72
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = "29/12/2025 04:57"
$ws.Range("C72").Value = "'955"
$ws.Range("D72").Value = "Inglês"
$ws.Range("E72").Value = "Verbs"
$ws.Range("F72").Value = 1

# Row 73
$ws.Range("A73").Value = 72
$ws.Range("B73").Value = "29/12/2025 04:58"
$ws.Range("C73").Value = "'954"
$ws.Range("D73").Value = "Inglês"
$ws.Range("E73").Value = "Pronouns"
$ws.Range("F73").Value = 1

# Row 74
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = "29/12/2025 05:00"
$ws.Range("C74").Value = "'953"
$ws.Range("D74").Value = "Inglês"
$ws.Range("E74").Value = "Semantic"
$ws.Range("F74").Value = 0

# Row 75
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = "29/12/2025 05:05"
$ws.Range("C75").Value = "'877"
$ws.Range("D75").Value = "Matemática Básica"
$ws.Range("E75").Value = "Média Ponderada"
$ws.Range("F75").Value = 1

# Row 76
$ws.Range("A76").Value = 75
$ws.Range("B76").Value = "29/12/2025 05:30"
$ws.Range("C76").Value = 1323
$ws.Range("D76").Value = "Matemática Financeira"
$ws.Range("E76").Value = "Equivalência de Capitais"
$ws.Range("F76").Value = 1

# Row 77
$ws.Range("A77").Value = 76
$ws.Range("B77").Value = "29/12/2025 05:35"
$ws.Range("C77").Value = 1330
$ws.Range("D77").Value = "Matemática Financeira"
$ws.Range("E77").Value = "Equivalência de Capitais"
$ws.Range("F77").Value = 0

# Row 78
$ws.Range("A78").Value = 77
$ws.Range("B78").Value = "29/12/2025 05:36"
$ws.Range("C78").Value = 1363
$ws.Range("D78").Value = "Matemática Financeira"
$ws.Range("E78").Value = "Análise de Investimentos"
$ws.Range("F78").Value = 0

# Row 79
$ws.Range("A79").Value = 78
$ws.Range("B79").Value = "29/12/2025 05:38"
$ws.Range("C79").Value = "'176"
$ws.Range("D79").Value = "Português"
$ws.Range("E79").Value = "Dois-Pontos"
$ws.Range("F79").Value = 1

# Row 80
$ws.Range("A80").Value = 79
$ws.Range("B80").Value = "29/12/2025 05:40"
$ws.Range("C80").Value = 240
$ws.Range("D80").Value = "Português"
$ws.Range("E80").Value = "Crase"
$ws.Range("F80").Value = 0

# Row 81
$ws.Range("A81").Value = 80
$ws.Range("B81").Value = "29/12/2025 05:41"
$ws.Range("C81").Value = "'200"
$ws.Range("D81").Value = "Português"
$ws.Range("E81").Value = "Concordância Verbal"
$ws.Range("F81").Value = 1

Write-Host "done"